$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 64 (shifts old rows 64-65 down to 66-67)
$ws.Rows.Item(64).Resize(2).Insert()

# Row 64 - new data (Primera)
$ws.Cells.Item(64, 1).Value = 9
$ws.Cells.Item(64, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(64, 3).Value = "Metropolitana"
$ws.Cells.Item(64, 4).Value = 44585
$ws.Cells.Item(64, 5).Value = 13
$ws.Cells.Item(64, 6).Value = 100114002
$ws.Cells.Item(64, 7).Value = "Camote"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 1330
$ws.Cells.Item(64, 11).Value = 11000
$ws.Cells.Item(64, 12).Value = 12000
$ws.Cells.Item(64, 13).Value = 11500
$ws.Cells.Item(64, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 639
$ws.Cells.Item(64, 17).Value = 18
$ws.Cells.Item(64, 18).Value = "Hortaliza"

# Row 65 - new data (Segunda)
$ws.Cells.Item(65, 1).Value = 9
$ws.Cells.Item(65, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(65, 3).Value = "Metropolitana"
$ws.Cells.Item(65, 4).Value = 44585
$ws.Cells.Item(65, 5).Value = 13
$ws.Cells.Item(65, 6).Value = 100114002
$ws.Cells.Item(65, 7).Value = "Camote"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Segunda"
$ws.Cells.Item(65, 10).Value = 610
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = 10000
$ws.Cells.Item(65, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(65, 15).Value = "Perú"
$ws.Cells.Item(65, 16).Value = 556
$ws.Cells.Item(65, 17).Value = 18
$ws.Cells.Item(65, 18).Value = "Hortaliza"
